# New weekly Apio / Vega Modelo de Temuco price observation was added to
# the source data. In the sheet this lands as a new data row inserted at
# row 465 (row 1 is the header), which pushes every existing row from 465
# down to 466 and extends the last row (old 524) to 525.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 465, shifting rows 465:524 down
# to 466:525 (matches the new <dimension ref="A1:R525"/>).
$ws.Rows.Item(465).EntireRow.Insert()

# Populate the newly inserted row 465 with the new record's data.
$ws.Cells.Item(465, 1).Value = 10
$ws.Cells.Item(465, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(465, 3).Value = "La Araucanía"
$ws.Cells.Item(465, 4).Value = 45142
$ws.Cells.Item(465, 5).Value = 9
$ws.Cells.Item(465, 6).Value = 100112017
$ws.Cells.Item(465, 7).Value = "Apio"
$ws.Cells.Item(465, 8).Value = "Americana (o)"
$ws.Cells.Item(465, 9).Value = "Primera"
$ws.Cells.Item(465, 10).Value = 40
$ws.Cells.Item(465, 11).Value = 7000
$ws.Cells.Item(465, 12).Value = 7000
$ws.Cells.Item(465, 13).Value = 7000
$ws.Cells.Item(465, 14).Value = "`$/docena de matas"
$ws.Cells.Item(465, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(465, 16).Value = 1167
$ws.Cells.Item(465, 17).Value = 6
$ws.Cells.Item(465, 18).Value = "Hortaliza"
